$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.299.81'
$ws.Range("E2").Value = '  +2.58%  '
$ws.Range("D3").Value = '1.803.03'
$ws.Range("E3").Value = '  +3.59%  '
$ws.Range("E4").Value = '  -0.83%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '340.14'
$ws.Range("E5").Value = '  +2.54%  '
$ws.Range("E6").Value = '  -0.26%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4591'
$ws.Range("E7").Value = '  +19.86%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3821'
$ws.Range("E8").Value = '  +13.97%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '45.20'
$ws.Range("E9").Value = '  -0.61%  '
$ws.Range("E10").Value = '  +5.47%  '
$ws.Range("E11").Value = '  +6.43%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '22.58'
$ws.Range("E12").Value = '  +2.26%  '
$ws.Range("E13").Value = '  -0.83%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.381'
$ws.Range("E14").Value = '  +4.33%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.557'
$ws.Range("E15").Value = '  +7.62%  '
$ws.Range("D16").Value = '1.804.79'
$ws.Range("E16").Value = '  +3.09%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001096'
$ws.Range("E17").Value = '  +4.49%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.06738'
$ws.Range("E18").Value = '  +2.21%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '81.70'
$ws.Range("E19").Value = '  +3.68%  '
$ws.Range("E20").Value = '  -0.36%  '
$ws.Range("E21").Value = '  +5.26%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.440'
$ws.Range("E22").Value = '  +4.89%  '
$ws.Range("D23").Value = '28.297.38'
$ws.Range("E23").Value = '  +2.16%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.93'
$ws.Range("E24").Value = '  +3.73%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.420'
$ws.Range("E25").Value = '  +0.36%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '20.78'
$ws.Range("E26").Value = '  +5.76%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '153.42'
$ws.Range("E27").Value = '  -0.26%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.386'
$ws.Range("E28").Value = '  +4.95%  '
$ws.Range("D29").Value = '2.011.36'
$ws.Range("E29").Value = '  +3.20%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '133.26'
$ws.Range("E30").Value = '  +2.60%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.253'
$ws.Range("E31").Value = '  -0.81%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.033'
$ws.Range("E32").Value = '  +0.72%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.09541'
$ws.Range("E33").Value = '  +10.14%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.880'
$ws.Range("E34").Value = '  +2.38%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.2310'
$ws.Range("E35").Value = '  +11.10%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '12.16'
$ws.Range("E36").Value = '  +1.75%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.302'
$ws.Range("E37").Value = '  +4.32%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02361'
$ws.Range("E38").Value = '  +5.02%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.06364'
$ws.Range("E39").Value = '  +5.42%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.6659'
$ws.Range("E40").Value = '  +3.39%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.241'
$ws.Range("E41").Value = '  +3.75%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '8.436'
$ws.Range("E42").Value = '  +6.48%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.485'
$ws.Range("E43").Value = '  -3.07%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '14.25'
$ws.Range("E44").Value = '  +5.16%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.9997'
$ws.Range("E45").Value = '  -0.20%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.6153'
$ws.Range("E46").Value = '  +3.30%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.878'
$ws.Range("E47").Value = '  +1.96%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '130.98'
$ws.Range("E48").Value = '  +3.89%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.047'
$ws.Range("E49").Value = '  +3.89%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.07188'
$ws.Range("E50").Value = '  +3.54%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.182'
$ws.Range("E51").Value = '  +3.14%  '
